# =====================================================================================
# Lithuania A Lyga - base update (commit: "Atualizacao de bases das ligas, do dia: 17-03-2024 as 10:24")
# =====================================================================================
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -------------------------------------------------------------------------------------
# 1) Two fixtures (id 24/25 and id 100/101/102) had been written to the wrong rows.
#    Re-key rows 26<->27 and rotate rows 102->103->104->102 so each "id" (col A, untouched)
#    lines up with the correct match data.
# -------------------------------------------------------------------------------------

# Row 26
$ws.Cells.Item(26, 2).Value2 = 6732711
$ws.Cells.Item(26, 3).Value2 = "Lithuania A Lyga"
$ws.Cells.Item(26, 4).Value2 = "Lithuania A Lyga"
$ws.Cells.Item(26, 5).Value2 = 45109.58333333334
$ws.Cells.Item(26, 6).Value2 = "Banga Gargzdai"
$ws.Cells.Item(26, 7).Value2 = "FK Zalgiris Vilnius"
$ws.Cells.Item(26, 8).Value2 = 1
$ws.Cells.Item(26, 9).Value2 = 4
$ws.Cells.Item(26, 10).Value2 = "A"
$ws.Cells.Item(26, 11).Value2 = 5
$ws.Cells.Item(26, 12).Value2 = 3.6
$ws.Cells.Item(26, 13).Value2 = 1.571
$ws.Cells.Item(26, 14).Value2 = 11
$ws.Cells.Item(26, 15).Value2 = 4.75
$ws.Cells.Item(26, 16).Value2 = 1.25
$ws.Cells.Item(26, 17).Value2 = 1.5
$ws.Cells.Item(26, 18).Value2 = 1.975
$ws.Cells.Item(26, 19).Value2 = 1.825
$ws.Cells.Item(26, 20).Value2 = 2.5
$ws.Cells.Item(26, 21).Value2 = 1.8
$ws.Cells.Item(26, 22).Value2 = 2
$ws.Cells.Item(26, 23).Value2 = -1
$ws.Cells.Item(26, 24).Value2 = -1
$ws.Cells.Item(26, 25).Value2 = 0.25
$ws.Cells.Item(26, 26).Value2 = -1
$ws.Cells.Item(26, 27).Value2 = 0.825
$ws.Cells.Item(26, 28).Value2 = 0.8
$ws.Cells.Item(26, 29).Value2 = -1

# Row 27
$ws.Cells.Item(27, 2).Value2 = 6732773
$ws.Cells.Item(27, 3).Value2 = "Lithuania A Lyga"
$ws.Cells.Item(27, 4).Value2 = "Lithuania A Lyga"
$ws.Cells.Item(27, 5).Value2 = 45109.58333333334
$ws.Cells.Item(27, 6).Value2 = "Suduva Marijampole"
$ws.Cells.Item(27, 7).Value2 = "Hegelmann Litauen"
$ws.Cells.Item(27, 8).Value2 = 0
$ws.Cells.Item(27, 9).Value2 = 1
$ws.Cells.Item(27, 10).Value2 = "A"
$ws.Cells.Item(27, 11).Value2 = 5
$ws.Cells.Item(27, 12).Value2 = 3.8
$ws.Cells.Item(27, 13).Value2 = 1.533
$ws.Cells.Item(27, 14).Value2 = 5
$ws.Cells.Item(27, 15).Value2 = 4.2
$ws.Cells.Item(27, 16).Value2 = 1.533
$ws.Cells.Item(27, 17).Value2 = 1
$ws.Cells.Item(27, 18).Value2 = 1.875
$ws.Cells.Item(27, 19).Value2 = 1.925
$ws.Cells.Item(27, 20).Value2 = 2.5
$ws.Cells.Item(27, 21).Value2 = 1.9
$ws.Cells.Item(27, 22).Value2 = 1.9
$ws.Cells.Item(27, 23).Value2 = -1
$ws.Cells.Item(27, 24).Value2 = -1
$ws.Cells.Item(27, 25).Value2 = 0.5329999999999999
$ws.Cells.Item(27, 26).Value2 = 0
$ws.Cells.Item(27, 27).Value2 = -0
$ws.Cells.Item(27, 28).Value2 = -1
$ws.Cells.Item(27, 29).Value2 = 0.8999999999999999

# Row 102
$ws.Cells.Item(102, 2).Value2 = 6732834
$ws.Cells.Item(102, 3).Value2 = "Lithuania A Lyga"
$ws.Cells.Item(102, 4).Value2 = "Lithuania A Lyga"
$ws.Cells.Item(102, 5).Value2 = 45242.41319444445
$ws.Cells.Item(102, 6).Value2 = "Panevezys"
$ws.Cells.Item(102, 7).Value2 = "FK Dziugas Telsiai"
$ws.Cells.Item(102, 8).Value2 = 0
$ws.Cells.Item(102, 9).Value2 = 0
$ws.Cells.Item(102, 10).Value2 = "D"
$ws.Cells.Item(102, 11).Value2 = 1.25
$ws.Cells.Item(102, 12).Value2 = 5.5
$ws.Cells.Item(102, 13).Value2 = 7.5
$ws.Cells.Item(102, 14).Value2 = 1.45
$ws.Cells.Item(102, 15).Value2 = 4.5
$ws.Cells.Item(102, 16).Value2 = 5
$ws.Cells.Item(102, 17).Value2 = -1
$ws.Cells.Item(102, 18).Value2 = 1.775
$ws.Cells.Item(102, 19).Value2 = 2.025
$ws.Cells.Item(102, 20).Value2 = 2.5
$ws.Cells.Item(102, 21).Value2 = 1.875
$ws.Cells.Item(102, 22).Value2 = 1.925
$ws.Cells.Item(102, 23).Value2 = -1
$ws.Cells.Item(102, 24).Value2 = 3.5
$ws.Cells.Item(102, 25).Value2 = -1
$ws.Cells.Item(102, 26).Value2 = -1
$ws.Cells.Item(102, 27).Value2 = 1.025
$ws.Cells.Item(102, 28).Value2 = -1
$ws.Cells.Item(102, 29).Value2 = 0.925

# Row 103
$ws.Cells.Item(103, 2).Value2 = 7465686
$ws.Cells.Item(103, 3).Value2 = "Lithuania A Lyga"
$ws.Cells.Item(103, 4).Value2 = "Lithuania A Lyga"
$ws.Cells.Item(103, 5).Value2 = 45242.41319444445
$ws.Cells.Item(103, 6).Value2 = "FK Kauno Zalgiris"
$ws.Cells.Item(103, 7).Value2 = "Hegelmann Litauen"
$ws.Cells.Item(103, 8).Value2 = 4
$ws.Cells.Item(103, 9).Value2 = 2
$ws.Cells.Item(103, 10).Value2 = "H"
$ws.Cells.Item(103, 11).Value2 = 2.3
$ws.Cells.Item(103, 12).Value2 = 4
$ws.Cells.Item(103, 13).Value2 = 2.3
$ws.Cells.Item(103, 14).Value2 = 2.55
$ws.Cells.Item(103, 15).Value2 = 4
$ws.Cells.Item(103, 16).Value2 = 2.2
$ws.Cells.Item(103, 17).Value2 = 0.25
$ws.Cells.Item(103, 18).Value2 = 1.8
$ws.Cells.Item(103, 19).Value2 = 2
$ws.Cells.Item(103, 20).Value2 = 2.75
$ws.Cells.Item(103, 21).Value2 = 1.85
$ws.Cells.Item(103, 22).Value2 = 1.95
$ws.Cells.Item(103, 23).Value2 = 1.55
$ws.Cells.Item(103, 24).Value2 = -1
$ws.Cells.Item(103, 25).Value2 = -1
$ws.Cells.Item(103, 26).Value2 = 0.8
$ws.Cells.Item(103, 27).Value2 = -1
$ws.Cells.Item(103, 28).Value2 = 0.8500000000000001
$ws.Cells.Item(103, 29).Value2 = -1

# Row 104
$ws.Cells.Item(104, 2).Value2 = 6732836
$ws.Cells.Item(104, 3).Value2 = "Lithuania A Lyga"
$ws.Cells.Item(104, 4).Value2 = "Lithuania A Lyga"
$ws.Cells.Item(104, 5).Value2 = 45242.41319444445
$ws.Cells.Item(104, 6).Value2 = "FK Siauliai"
$ws.Cells.Item(104, 7).Value2 = "Banga Gargzdai"
$ws.Cells.Item(104, 8).Value2 = 3
$ws.Cells.Item(104, 9).Value2 = 0
$ws.Cells.Item(104, 10).Value2 = "H"
$ws.Cells.Item(104, 11).Value2 = 1.222
$ws.Cells.Item(104, 12).Value2 = 5.5
$ws.Cells.Item(104, 13).Value2 = 9
$ws.Cells.Item(104, 14).Value2 = 1.363
$ws.Cells.Item(104, 15).Value2 = 4.5
$ws.Cells.Item(104, 16).Value2 = 7
$ws.Cells.Item(104, 17).Value2 = -1.25
$ws.Cells.Item(104, 18).Value2 = 1.9
$ws.Cells.Item(104, 19).Value2 = 1.9
$ws.Cells.Item(104, 20).Value2 = 2.5
$ws.Cells.Item(104, 21).Value2 = 1.975
$ws.Cells.Item(104, 22).Value2 = 1.825
$ws.Cells.Item(104, 23).Value2 = 0.363
$ws.Cells.Item(104, 24).Value2 = -1
$ws.Cells.Item(104, 25).Value2 = -1
$ws.Cells.Item(104, 26).Value2 = 0.8999999999999999
$ws.Cells.Item(104, 27).Value2 = -1
$ws.Cells.Item(104, 28).Value2 = 0.9750000000000001
$ws.Cells.Item(104, 29).Value2 = -1

# -------------------------------------------------------------------------------------
# 2) Row 115 (id 113) was missing its full-time score / result / closing-odds-movement
#    columns (H, I, J, AB, AC) and had placeholder zeros in N:AA - fill in the final data.
# -------------------------------------------------------------------------------------

# Row 115
$ws.Cells.Item(115, 2).Value2 = 7862909
$ws.Cells.Item(115, 3).Value2 = "Lithuania A Lyga"
$ws.Cells.Item(115, 4).Value2 = "Lithuania A Lyga"
$ws.Cells.Item(115, 5).Value2 = 45363.54166666666
$ws.Cells.Item(115, 6).Value2 = "FK Kauno Zalgiris"
$ws.Cells.Item(115, 7).Value2 = "FK Transinvest"
$ws.Cells.Item(115, 8).Value2 = 2
$ws.Cells.Item(115, 9).Value2 = 1
$ws.Cells.Item(115, 10).Value2 = "H"
$ws.Cells.Item(115, 11).Value2 = 2
$ws.Cells.Item(115, 12).Value2 = 3.25
$ws.Cells.Item(115, 13).Value2 = 3.25
$ws.Cells.Item(115, 14).Value2 = 1.571
$ws.Cells.Item(115, 15).Value2 = 3.6
$ws.Cells.Item(115, 16).Value2 = 4.75
$ws.Cells.Item(115, 17).Value2 = -0.75
$ws.Cells.Item(115, 18).Value2 = 1.825
$ws.Cells.Item(115, 19).Value2 = 1.975
$ws.Cells.Item(115, 20).Value2 = 2.5
$ws.Cells.Item(115, 21).Value2 = 1.825
$ws.Cells.Item(115, 22).Value2 = 1.975
$ws.Cells.Item(115, 23).Value2 = 0.571
$ws.Cells.Item(115, 24).Value2 = -1
$ws.Cells.Item(115, 25).Value2 = -1
$ws.Cells.Item(115, 26).Value2 = 0.4125
$ws.Cells.Item(115, 27).Value2 = -0.5
$ws.Cells.Item(115, 28).Value2 = 0.825
$ws.Cells.Item(115, 29).Value2 = -1

# -------------------------------------------------------------------------------------
# 3) Append 8 new fixtures as rows 116:123 (ids 114-121). Apply the same visual style used
#    by every other data row first (col A: bold + thin border + centered/top; col E: the
#    custom date/time number format) by copying it from row 2, which reuses the existing
#    style indices in styles.xml instead of registering new, unused ones.
# -------------------------------------------------------------------------------------

$ws.Range("A2").Copy()
$ws.Range("A116:A123").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E116:E123").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 116
$ws.Cells.Item(116, 1).Value2 = 114
$ws.Cells.Item(116, 2).Value2 = 7862910
$ws.Cells.Item(116, 3).Value2 = "Lithuania A Lyga"
$ws.Cells.Item(116, 4).Value2 = "Lithuania A Lyga"
$ws.Cells.Item(116, 5).Value2 = 45364.5
$ws.Cells.Item(116, 6).Value2 = "FK Dainava Alytus"
$ws.Cells.Item(116, 7).Value2 = "FK Dziugas Telsiai"
$ws.Cells.Item(116, 8).Value2 = 0
$ws.Cells.Item(116, 9).Value2 = 0
$ws.Cells.Item(116, 10).Value2 = "D"
$ws.Cells.Item(116, 11).Value2 = 2.3
$ws.Cells.Item(116, 12).Value2 = 2.7
$ws.Cells.Item(116, 13).Value2 = 3.25
$ws.Cells.Item(116, 14).Value2 = 2.4
$ws.Cells.Item(116, 15).Value2 = 2.75
$ws.Cells.Item(116, 16).Value2 = 3
$ws.Cells.Item(116, 17).Value2 = -0.25
$ws.Cells.Item(116, 18).Value2 = 2.075
$ws.Cells.Item(116, 19).Value2 = 1.725
$ws.Cells.Item(116, 20).Value2 = 1.75
$ws.Cells.Item(116, 21).Value2 = 1.775
$ws.Cells.Item(116, 22).Value2 = 2.025
$ws.Cells.Item(116, 23).Value2 = -1
$ws.Cells.Item(116, 24).Value2 = 1.75
$ws.Cells.Item(116, 25).Value2 = -1
$ws.Cells.Item(116, 26).Value2 = -0.5
$ws.Cells.Item(116, 27).Value2 = 0.3625
$ws.Cells.Item(116, 28).Value2 = -1
$ws.Cells.Item(116, 29).Value2 = 1.025

# Row 117
$ws.Cells.Item(117, 1).Value2 = 115
$ws.Cells.Item(117, 2).Value2 = 7862911
$ws.Cells.Item(117, 3).Value2 = "Lithuania A Lyga"
$ws.Cells.Item(117, 4).Value2 = "Lithuania A Lyga"
$ws.Cells.Item(117, 5).Value2 = 45364.54166666666
$ws.Cells.Item(117, 6).Value2 = "Hegelmann Litauen"
$ws.Cells.Item(117, 7).Value2 = "FK Siauliai"
$ws.Cells.Item(117, 8).Value2 = 2
$ws.Cells.Item(117, 9).Value2 = 2
$ws.Cells.Item(117, 10).Value2 = "D"
$ws.Cells.Item(117, 11).Value2 = 2.15
$ws.Cells.Item(117, 12).Value2 = 3.1
$ws.Cells.Item(117, 13).Value2 = 3.1
$ws.Cells.Item(117, 14).Value2 = 2.45
$ws.Cells.Item(117, 15).Value2 = 2.9
$ws.Cells.Item(117, 16).Value2 = 3
$ws.Cells.Item(117, 17).Value2 = 0
$ws.Cells.Item(117, 18).Value2 = 1.725
$ws.Cells.Item(117, 19).Value2 = 2.075
$ws.Cells.Item(117, 20).Value2 = 2.5
$ws.Cells.Item(117, 21).Value2 = 2.025
$ws.Cells.Item(117, 22).Value2 = 1.775
$ws.Cells.Item(117, 23).Value2 = -1
$ws.Cells.Item(117, 24).Value2 = 1.9
$ws.Cells.Item(117, 25).Value2 = -1
$ws.Cells.Item(117, 26).Value2 = 0
$ws.Cells.Item(117, 27).Value2 = -0
$ws.Cells.Item(117, 28).Value2 = 1.025
$ws.Cells.Item(117, 29).Value2 = -1

# Row 118
$ws.Cells.Item(118, 1).Value2 = 116
$ws.Cells.Item(118, 2).Value2 = 7862036
$ws.Cells.Item(118, 3).Value2 = "Lithuania A Lyga"
$ws.Cells.Item(118, 4).Value2 = "Lithuania A Lyga"
$ws.Cells.Item(118, 5).Value2 = 45364.54166666666
$ws.Cells.Item(118, 6).Value2 = "Banga Gargzdai"
$ws.Cells.Item(118, 7).Value2 = "FK Zalgiris Vilnius"
$ws.Cells.Item(118, 8).Value2 = 1
$ws.Cells.Item(118, 9).Value2 = 4
$ws.Cells.Item(118, 10).Value2 = "A"
$ws.Cells.Item(118, 11).Value2 = 8
$ws.Cells.Item(118, 12).Value2 = 4.5
$ws.Cells.Item(118, 13).Value2 = 1.3
$ws.Cells.Item(118, 14).Value2 = 6.5
$ws.Cells.Item(118, 15).Value2 = 4.5
$ws.Cells.Item(118, 16).Value2 = 1.333
$ws.Cells.Item(118, 17).Value2 = 1.25
$ws.Cells.Item(118, 18).Value2 = 2
$ws.Cells.Item(118, 19).Value2 = 1.8
$ws.Cells.Item(118, 20).Value2 = 2.5
$ws.Cells.Item(118, 21).Value2 = 1.825
$ws.Cells.Item(118, 22).Value2 = 1.975
$ws.Cells.Item(118, 23).Value2 = -1
$ws.Cells.Item(118, 24).Value2 = -1
$ws.Cells.Item(118, 25).Value2 = 0.333
$ws.Cells.Item(118, 26).Value2 = -1
$ws.Cells.Item(118, 27).Value2 = 0.8
$ws.Cells.Item(118, 28).Value2 = 0.825
$ws.Cells.Item(118, 29).Value2 = -1

# Row 119
$ws.Cells.Item(119, 1).Value2 = 117
$ws.Cells.Item(119, 2).Value2 = 7862037
$ws.Cells.Item(119, 3).Value2 = "Lithuania A Lyga"
$ws.Cells.Item(119, 4).Value2 = "Lithuania A Lyga"
$ws.Cells.Item(119, 5).Value2 = 45364.63541666666
$ws.Cells.Item(119, 6).Value2 = "Suduva Marijampole"
$ws.Cells.Item(119, 7).Value2 = "Panevezys"
$ws.Cells.Item(119, 8).Value2 = 1
$ws.Cells.Item(119, 9).Value2 = 0
$ws.Cells.Item(119, 10).Value2 = "H"
$ws.Cells.Item(119, 11).Value2 = 7.5
$ws.Cells.Item(119, 12).Value2 = 4
$ws.Cells.Item(119, 13).Value2 = 1.363
$ws.Cells.Item(119, 14).Value2 = 3.1
$ws.Cells.Item(119, 15).Value2 = 3.25
$ws.Cells.Item(119, 16).Value2 = 2.1
$ws.Cells.Item(119, 17).Value2 = 0.25
$ws.Cells.Item(119, 18).Value2 = 1.875
$ws.Cells.Item(119, 19).Value2 = 1.925
$ws.Cells.Item(119, 20).Value2 = 2
$ws.Cells.Item(119, 21).Value2 = 1.95
$ws.Cells.Item(119, 22).Value2 = 1.85
$ws.Cells.Item(119, 23).Value2 = 2.1
$ws.Cells.Item(119, 24).Value2 = -1
$ws.Cells.Item(119, 25).Value2 = -1
$ws.Cells.Item(119, 26).Value2 = 0.875
$ws.Cells.Item(119, 27).Value2 = -1
$ws.Cells.Item(119, 28).Value2 = -1
$ws.Cells.Item(119, 29).Value2 = 0.8500000000000001

# Row 120
$ws.Cells.Item(120, 1).Value2 = 118
$ws.Cells.Item(120, 2).Value2 = 7862038
$ws.Cells.Item(120, 3).Value2 = "Lithuania A Lyga"
$ws.Cells.Item(120, 4).Value2 = "Lithuania A Lyga"
$ws.Cells.Item(120, 5).Value2 = 45367.375
$ws.Cells.Item(120, 6).Value2 = "FK Dziugas Telsiai"
$ws.Cells.Item(120, 7).Value2 = "Suduva Marijampole"
$ws.Cells.Item(120, 8).Value2 = 1
$ws.Cells.Item(120, 9).Value2 = 0
$ws.Cells.Item(120, 10).Value2 = "H"
$ws.Cells.Item(120, 11).Value2 = 2.5
$ws.Cells.Item(120, 12).Value2 = 3.2
$ws.Cells.Item(120, 13).Value2 = 2.5
$ws.Cells.Item(120, 14).Value2 = 3.2
$ws.Cells.Item(120, 15).Value2 = 3.1
$ws.Cells.Item(120, 16).Value2 = 2.05
$ws.Cells.Item(120, 17).Value2 = 0.25
$ws.Cells.Item(120, 18).Value2 = 1.95
$ws.Cells.Item(120, 19).Value2 = 1.85
$ws.Cells.Item(120, 20).Value2 = 2
$ws.Cells.Item(120, 21).Value2 = 1.9
$ws.Cells.Item(120, 22).Value2 = 1.9
$ws.Cells.Item(120, 23).Value2 = 2.2
$ws.Cells.Item(120, 24).Value2 = -1
$ws.Cells.Item(120, 25).Value2 = -1
$ws.Cells.Item(120, 26).Value2 = 0.95
$ws.Cells.Item(120, 27).Value2 = -1
$ws.Cells.Item(120, 28).Value2 = -1
$ws.Cells.Item(120, 29).Value2 = 0.8999999999999999

# Row 121
$ws.Cells.Item(121, 1).Value2 = 119
$ws.Cells.Item(121, 2).Value2 = 7862912
$ws.Cells.Item(121, 3).Value2 = "Lithuania A Lyga"
$ws.Cells.Item(121, 4).Value2 = "Lithuania A Lyga"
$ws.Cells.Item(121, 5).Value2 = 45367.4375
$ws.Cells.Item(121, 6).Value2 = "Panevezys"
$ws.Cells.Item(121, 7).Value2 = "Banga Gargzdai"
$ws.Cells.Item(121, 8).Value2 = 0
$ws.Cells.Item(121, 9).Value2 = 1
$ws.Cells.Item(121, 10).Value2 = "A"
$ws.Cells.Item(121, 11).Value2 = 1.2
$ws.Cells.Item(121, 12).Value2 = 6.5
$ws.Cells.Item(121, 13).Value2 = 8
$ws.Cells.Item(121, 14).Value2 = 1.4
$ws.Cells.Item(121, 15).Value2 = 5
$ws.Cells.Item(121, 16).Value2 = 5
$ws.Cells.Item(121, 17).Value2 = -1
$ws.Cells.Item(121, 18).Value2 = 1.775
$ws.Cells.Item(121, 19).Value2 = 2.025
$ws.Cells.Item(121, 20).Value2 = 2.25
$ws.Cells.Item(121, 21).Value2 = 1.95
$ws.Cells.Item(121, 22).Value2 = 1.85
$ws.Cells.Item(121, 23).Value2 = -1
$ws.Cells.Item(121, 24).Value2 = -1
$ws.Cells.Item(121, 25).Value2 = 4
$ws.Cells.Item(121, 26).Value2 = -1
$ws.Cells.Item(121, 27).Value2 = 1.025
$ws.Cells.Item(121, 28).Value2 = -1
$ws.Cells.Item(121, 29).Value2 = 0.8500000000000001

# Row 122
$ws.Cells.Item(122, 1).Value2 = 120
$ws.Cells.Item(122, 2).Value2 = 7862914
$ws.Cells.Item(122, 3).Value2 = "Lithuania A Lyga"
$ws.Cells.Item(122, 4).Value2 = "Lithuania A Lyga"
$ws.Cells.Item(122, 5).Value2 = 45368.41666666666
$ws.Cells.Item(122, 6).Value2 = "FK Siauliai"
$ws.Cells.Item(122, 7).Value2 = "FK Kauno Zalgiris"
$ws.Cells.Item(122, 11).Value2 = 2.5
$ws.Cells.Item(122, 12).Value2 = 2.875
$ws.Cells.Item(122, 13).Value2 = 2.75
$ws.Cells.Item(122, 14).Value2 = 2.3
$ws.Cells.Item(122, 15).Value2 = 3
$ws.Cells.Item(122, 16).Value2 = 3.1
$ws.Cells.Item(122, 17).Value2 = -0.25
$ws.Cells.Item(122, 18).Value2 = 2.025
$ws.Cells.Item(122, 19).Value2 = 1.775
$ws.Cells.Item(122, 20).Value2 = 2.25
$ws.Cells.Item(122, 21).Value2 = 2
$ws.Cells.Item(122, 22).Value2 = 1.8
$ws.Cells.Item(122, 23).Value2 = 0
$ws.Cells.Item(122, 24).Value2 = 0
$ws.Cells.Item(122, 25).Value2 = 0
$ws.Cells.Item(122, 26).Value2 = 0
$ws.Cells.Item(122, 27).Value2 = 0

# Row 123
$ws.Cells.Item(123, 1).Value2 = 121
$ws.Cells.Item(123, 2).Value2 = 7862039
$ws.Cells.Item(123, 3).Value2 = "Lithuania A Lyga"
$ws.Cells.Item(123, 4).Value2 = "Lithuania A Lyga"
$ws.Cells.Item(123, 5).Value2 = 45368.5625
$ws.Cells.Item(123, 6).Value2 = "FK Zalgiris Vilnius"
$ws.Cells.Item(123, 7).Value2 = "Hegelmann Litauen"
$ws.Cells.Item(123, 11).Value2 = 1.615
$ws.Cells.Item(123, 12).Value2 = 3.8
$ws.Cells.Item(123, 13).Value2 = 4.333
$ws.Cells.Item(123, 14).Value2 = 1.615
$ws.Cells.Item(123, 15).Value2 = 3.75
$ws.Cells.Item(123, 16).Value2 = 4.2
$ws.Cells.Item(123, 17).Value2 = -0.75
$ws.Cells.Item(123, 18).Value2 = 1.85
$ws.Cells.Item(123, 19).Value2 = 1.95
$ws.Cells.Item(123, 20).Value2 = 2.75
$ws.Cells.Item(123, 21).Value2 = 1.975
$ws.Cells.Item(123, 22).Value2 = 1.825
$ws.Cells.Item(123, 23).Value2 = 0
$ws.Cells.Item(123, 24).Value2 = 0
$ws.Cells.Item(123, 25).Value2 = 0
$ws.Cells.Item(123, 26).Value2 = 0
$ws.Cells.Item(123, 27).Value2 = 0

